# Insert a new row at position 106, shifting existing rows 106-128 down to 107-129
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with its data
$ws.Cells.Item(106, 1).Value = 7
$ws.Cells.Item(106, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(106, 3).Value = "Ñuble"
$ws.Cells.Item(106, 4).Value = 44476
$ws.Cells.Item(106, 5).Value = 16
$ws.Cells.Item(106, 6).Value = 100112006
$ws.Cells.Item(106, 7).Value = "Repollo"
$ws.Cells.Item(106, 8).Value = "Crespo record"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 400
$ws.Cells.Item(106, 11).Value = 600
$ws.Cells.Item(106, 12).Value = 650
$ws.Cells.Item(106, 13).Value = 625
$ws.Cells.Item(106, 14).Value = "$/unidad"
$ws.Cells.Item(106, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(106, 16).Value = 625
$ws.Cells.Item(106, 17).Value = 1
$ws.Cells.Item(106, 18).Value = "Hortaliza"
